# Update stats for 2025-12 (row 25) in the iserv_stats workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6477
$ws.Range("C25").Value = 1008
$ws.Range("D25").Value = 6022492
$ws.Range("E25").Value = 929.8273892234059
$ws.Range("F25").Value = 9.947377355287724
$ws.Range("G25").Value = 7.462686567164178
$ws.Range("H25").Value = 26.12955475422662
